# Adds a new "Italy" worksheet (test data for Italy market), mirroring the
# existing "Germany" sheet layout (same columns / repeater rows), and
# updates the two market-specific cells. Also reproduces the selection /
# active-tab bookkeeping changes that Excel records as a side effect of
# this edit.

$wb = $excel.ActiveWorkbook

# --- 1. Record Germany's new selection (whole-sheet select) ---------------
$germany = $wb.Worksheets.Item("Germany")
$germany.Activate()
$germany.Cells.Select()

# --- 2. Build the new "Italy" sheet from a copy of "Germany" --------------
# Germany already has the full repeater list (incl. P32AR/P32DR) and the
# column widths that the new Italy sheet uses, so copying it and only
# touching the two market-specific cells reproduces the target layout.
$slovakia = $wb.Worksheets.Item("Slovakia")
$germany.Copy($null, $slovakia)

$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

$italy.Range("B2").Value = "Italy Market"

# The new part number cell loses its inherited "Color Codes" style in the
# authored change (no s= attribute on B4), so clear formatting first.
$italy.Range("B4").ClearFormats()
$italy.Range("B4").Value = "NGC-3145/T2156/T2158 "

# --- 3. Slovakia keeps its data but is no longer the selected tab; its
#        selection moves to B20 -------------------------------------------
$slovakia.Activate()
$slovakia.Range("B20").Select()

# --- 4. Italy becomes the active / selected tab, with B4 selected ---------
$italy.Activate()
$italy.Range("B4").Select()
